# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Bad Drivers table - update Critical Minutes and Good Roaming Calculation (%)
$ws.Range("C3").Value = 1214
$ws.Range("D3").Value = 98

# Totals row - Critical Minutes total
$ws.Range("C4").Value = 1214

# Good Drivers table - update Total Samples
$ws.Range("B12").Value = 449371
$ws.Range("B13").Value = 77999
